$d = $word.ActiveDocument

# The requested change trims the document-wide default run/paragraph
# formatting (w:docDefaults inside word/styles.xml) down to just a
# handful of properties. That block isn't reachable through the normal
# Style/Font/ParagraphFormat object model (those edit the "Normal"
# style's own w:rPr/w:pPr, not docDefaults), so we round-trip the whole
# package through Document.WordOpenXML and patch the docDefaults XML
# fragment directly.

$newDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$xml = $d.WordOpenXML

# Pull out the existing <w:docDefaults>...</w:docDefaults> element (attribute
# order/self-closing shorthand can vary slightly by serializer) and swap its
# whole contents for the trimmed version the diff asks for.
$pattern = '<w:docDefaults>.*?</w:docDefaults>'
$match = [System.Text.RegularExpressions.Regex]::Match($xml, $pattern, [System.Text.RegularExpressions.RegexOptions]::Singleline)

if ($match.Success) {
    $xml = $xml.Substring(0, $match.Index) + $newDefaults + $xml.Substring($match.Index + $match.Length)
    $d.WordOpenXML = $xml
    Write-Output "docDefaults updated"
} else {
    Write-Output "WARNING: docDefaults block not found"
}
